# Update odds values for rows 4, 6 and 9 (matches FlashScore on 2024-11-21)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 updates
$ws.Range("G4").Value = 5.1
$ws.Range("H4").Value = 3.45
$ws.Range("I4").Value = 1.65
$ws.Range("J4").Value = 5.1
$ws.Range("L4").Value = 2.22
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.34
$ws.Range("P4").Value = 2.75
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 1.65
$ws.Range("S4").Value = 1.39
$ws.Range("T4").Value = 2.55
$ws.Range("U4").Value = 1.93
$ws.Range("V4").Value = 1.7
$ws.Range("W4").Value = 12
$ws.Range("X4").Value = 29
$ws.Range("Y4").Value = 16.5
$ws.Range("AA4").Value = 55
$ws.Range("AB4").Value = 60
$ws.Range("AC4").Value = 8.5
$ws.Range("AD4").Value = 6.7
$ws.Range("AE4").Value = 17.5
$ws.Range("AF4").Value = 100
$ws.Range("AH4").Value = 7.1
$ws.Range("AM4").Value = 900
$ws.Range("AO4").Value = 29
$ws.Range("AP4").Value = 32
$ws.Range("AQ4").Value = 175
$ws.Range("AR4").Value = 200
$ws.Range("AS4").Value = 450
$ws.Range("AT4").Value = 2.52
$ws.Range("AU4").Value = 7.5
$ws.Range("AV4").Value = 70
$ws.Range("AX4").Value = 8
$ws.Range("AZ4").Value = 27

# Row 6 updates
$ws.Range("G6").Value = 3.2
$ws.Range("I6").Value = 2.2
$ws.Range("J6").Value = 3.75
$ws.Range("L6").Value = 2.88
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("U6").Value = 1.73
$ws.Range("V6").Value = 2
$ws.Range("W6").Value = 11
$ws.Range("X6").Value = 17
$ws.Range("Y6").Value = 12
$ws.Range("AA6").Value = 26
$ws.Range("AB6").Value = 34
$ws.Range("AH6").Value = 11
$ws.Range("AK6").Value = 17
$ws.Range("AM6").Value = 201
$ws.Range("AP6").Value = 26
$ws.Range("AR6").Value = 81
$ws.Range("AU6").Value = 8
$ws.Range("AX6").Value = 12

# Row 9 updates
$ws.Range("I9").Value = 3.25
$ws.Range("K9").Value = 2.2
$ws.Range("X9").Value = 10
$ws.Range("AB9").Value = 26
$ws.Range("AC9").Value = 11
$ws.Range("AK9").Value = 26
$ws.Range("AY9").Value = 26
